$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple value replacements on the first rows (indices stable so far) ---
$t.Cell(1,1).Range.Text = "0M"       # was 99.81
$t.Cell(2,1).Range.Text = "0M"       # was 0.28
$t.Cell(3,1).Range.Text = "0M"       # was 147
$t.Cell(4,1).Range.Text = "1300"     # was 601
$t.Cell(5,1).Range.Text = "0.00002"  # was 0.00003
$t.Cell(6,1).Range.Text = "0.00073"  # was 0.00070 (rows 7,8,9 get removed below)

# --- Step 2: remove the three now-redundant rows (0.00013, 0.00005, 0.00015) ---
# delete from highest index to lowest so earlier indices remain valid
$t.Rows.Item(9).Delete()
$t.Rows.Item(8).Delete()
$t.Rows.Item(7).Delete()

# After the deletions the table looks like (1-indexed):
#  6 -> 0.00073 (just updated above)
#  7 -> 0.00018 (unchanged)
#  8 -> 0.00021 (to become 0.00006)
#  9 -> 0.08230 (to become 0.00030)
# 10 -> 100.0 (unchanged)

$t.Cell(8,1).Range.Text = "0.00006"  # was 0.00021
$t.Cell(9,1).Range.Text = "0.00030"  # was 0.08230

# --- Step 3: insert three new rows right after row 9 (0.00030), before row 10 (100.0) ---
$r = $t.Rows.Add($t.Rows.Item(10))
$t.Cell(10,1).Range.Text = "0.00040"

$r = $t.Rows.Add($t.Rows.Item(11))
$t.Cell(11,1).Range.Text = "0.00051"

$r = $t.Rows.Add($t.Rows.Item(12))
$t.Cell(12,1).Range.Text = "0.27972"

# --- Step 4: collapse the three trailing multi-run summary rows into single values ---
# the row count is back to 46 (46 - 3 deleted + 3 added), so these are still at 44/45/46
$t.Cell(44,1).Range.Text = "99.81"   # was "350`t0.00020`t0.00073`t0.00042`t0.00013`t0.00030`t0.00040`t0.00051`t0.14709`t100.0"
$t.Cell(45,1).Range.Text = "0.28"    # was "99`t0.00002`t0.00008`t0.00006`t0.00001`t0.00005`t0.00007`t0.00007`t0.00579`t100.0"
$t.Cell(46,1).Range.Text = "147"     # was "250`t0.00008`t0.00032`t0.00018`t0.00006`t0.00013`t0.00017`t0.00022`t0.04454`t100.0"
